# Adds summary statistics (average/worst ratios) below the results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the |S*|/n column (J), bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Row 14: labeled summary metric - average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# Row 15: labeled summary metric - average of SC(S*)/SC(OPT)
# (re-use B14's formatting so no extra/duplicate cell styles are produced)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)

# Row 16: labeled summary metric - worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B14").Copy()
$ws.Range("B16").PasteSpecial(-4122)

# Row 17: labeled summary metric - worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$ws.Range("B14").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# Match the taller row height used for the bold 12pt summary rows
$ws.Range("A14:B17").RowHeight = 15.6

# Page setup tweaks
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on the newly added average cell
$ws.Range("J12").Select()
